$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Sheet2")

# Mark D3 (existing "foto1.jpg; foto2.jpg" answer) as wrap-text, and add two
# new example-answer cells below it for the photo-filename question.
$ws2.Range("D3").WrapText = $true
$ws2.Range("D4").Value = "foto23.jpg"
$ws2.Range("D5").Value = "foto2323.jpg;foto235.jpg"

# The question list ("vraag 2" .. "vraag 21") used to start at row 9 with a
# blank row 10 in the middle; it now starts right after the new D4/D5 answers
# at row 6, with no gap. Remove the stray blank row first, then close the
# remaining 3-row gap above the list.
$ws2.Rows.Item(10).Delete()
$ws2.Range("A6:A8").EntireRow.Delete()

# Keep the active selection sane on the now-shorter sheet.
$ws2.Range("D3").Select()
